# chore: adapt column header formatting to respective input file names
#
# 1) Rename the header cells in row 1 from the "_old"/"_new" suffix scheme
#    to the "_FV2310"/"_FV2404" (format-version based) suffix scheme.
# 2) Freeze the header row (row 1) so it stays visible while scrolling.
# 3) Turn the data range A1:U84 into a native Excel Table ("Table1") so the
#    header names double as the table's column headers / the range gets an
#    AutoFilter.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Rename header cells --------------------------------------------
$baseNames = @(
    "Segmentname",
    "Segmentgruppe",
    "Segment",
    "Datenelement",
    "Segment ID",
    "Code",
    "Qualifier",
    "Beschreibung",
    "Bedingungsausdruck",
    "Bedingung"
)

# Columns A..J hold the "_old" -> "_FV2310" (left / previous format version) headers
$oldCols = @("A", "B", "C", "D", "E", "F", "G", "H", "I", "J")
for ($i = 0; $i -lt $oldCols.Length; $i++) {
    $ws.Range($oldCols[$i] + "1").Value2 = $baseNames[$i] + "_FV2310"
}

# Column K is "diff" and stays untouched.

# Columns L..U hold the "_new" -> "_FV2404" (right / current format version) headers
$newCols = @("L", "M", "N", "O", "P", "Q", "R", "S", "T", "U")
for ($i = 0; $i -lt $newCols.Length; $i++) {
    $ws.Range($newCols[$i] + "1").Value2 = $baseNames[$i] + "_FV2404"
}

# --- 2) Freeze the header row -------------------------------------------
$ws.Activate() | Out-Null
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true

# --- 3) Convert the used range into an Excel Table ----------------------
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U84"), $null, 1)
$tbl.Name = "Table1"

Write-Host "Renamed headers, froze top row, and added Table1 over A1:U84"
